$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the raw pixel-coordinate values (trained YOLO produced new crops) ---
$ws.Range("C2").Value = 4961
$ws.Range("D2").Value = 4016
$ws.Range("C4").Value = 185
$ws.Range("D4").Value = 2205
$ws.Range("C5").Value = 1025
$ws.Range("D5").Value = 2754

# --- Apply a "Comma" (Millares) number style to the coordinate block C2:D5 ---
# Excel's built-in "Comma" style (xfId/cellStyleXf builtinId=3) is applied first;
# this is what creates the thousands separator + the associated style sheet
# records (numFmt 43, a duplicated font, a cellStyleXf and a cellStyle named
# "Comma"/"Millares").
$ws.Range("C2:D5").Style = "Comma"
$ws.Range("C2:D5").HorizontalAlignment = -4108
$ws.Range("C2:D5").NumberFormat = "_-* #,##0_-;\-* #,##0_-;_-* ""-""??_-;_-@_-"

# --- Apply the same family of format, with more decimals, to the ratio block C7:D10 ---
$ws.Range("C7:D10").Style = "Comma"
$ws.Range("C7:D10").HorizontalAlignment = -4108
$ws.Range("C7:D10").NumberFormat = "_-* #,##0.0000_-;\-* #,##0.0000_-;_-* ""-""??_-;_-@_-"

# --- Column widths for the newly-formatted coordinate columns ---
$ws.Range("C:D").ColumnWidth = 10.2

# --- Recalculate so the cached formula results match the new inputs ---
$excel.Calculate()

# --- Restore the selection the author ended up with after the edit ---
$ws.Range("G8").Select() | Out-Null
